# Repurpose geothermal as pumped hydro
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# ===========================================================================
# 1. Insert new "Pumped storage" worksheet before "SYC-SYEGC"
# ===========================================================================
$target = $wb.Worksheets.Item("SYC-SYEGC")
$ps = $wb.Worksheets.Add($target)
$ps.Name = "Pumped storage"

# ---- Titles (rows 1-3) ----------------------------------------------------
$ps.Range("A1").Value = "PUMPED STORAGE DEVELOPMENT IN INDIA"
$ps.Range("A2").Value = "(Installed Capacity above 25 MW)"
$ps.Range("A3").Value = "As on 30.09.2019"
$ps.Range("A1:A3").Font.Name = "Arial"
$ps.Range("A1:A3").Font.Size = 12

# ---- Table border (rows 4-19, cols A-E) ------------------------------------
$ps.Range("A4:E19").Borders.LineStyle = 1
$ps.Range("A4:E19").Borders.Weight = 2

# ---- Merged header "Installed Capacity" -----------------------------------
$ps.Range("C4:D4").Merge()
$ps.Range("C4").Value = "Installed Capacity"
$ps.Range("C4:D4").HorizontalAlignment = -4108
$ps.Range("C4:D4").VerticalAlignment = -4108
$ps.Range("C4:D4").Font.Bold = $true

# ---- Column headers (row 5) ------------------------------------------------
$ps.Range("A5").Value = "S.No."
$ps.Range("B5").Value = "SCHEMES"
$ps.Range("C5").Value = "no. of units `nX unit size (MW)"
$ps.Range("D5").Value = "MW"
$ps.Range("E5").Value = "Remarks"

$ps.Range("A5:B5").Font.Name = "Arial"
$ps.Range("A5:B5").Font.Size = 10
$ps.Range("A5:B5").Font.Bold = $true
$ps.Range("A5:B5").HorizontalAlignment = -4108

$ps.Range("C5").Font.Bold = $true
$ps.Range("C5").HorizontalAlignment = -4108
$ps.Range("C5").WrapText = $true

$ps.Range("D5:E5").Font.Bold = $true
$ps.Range("D5:E5").HorizontalAlignment = -4108

# ---- Section header rows 6-7 -----------------------------------------------
$ps.Range("A6").Value = "A. Schemes Constructed"
$ps.Range("A7").Value = "a) Working in pumping mode"
$ps.Range("A6:A7").Font.Name = "Arial"
$ps.Range("A6:A7").Font.Size = 12
$ps.Range("A6:A7").Font.Bold = $true

# ---- Data rows 8-13 (working schemes) --------------------------------------
$ps.Range("A8").Value = 1
$ps.Range("B8").Value = "Nagarjuna Sagar -Telangana"
$ps.Range("C8").Value = "7x100.80"
$ps.Range("D8").Value = 705.6

$ps.Range("A9").Value = 2
$ps.Range("B9").Value = "Kadamparai -T.N"
$ps.Range("C9").Value = "4x100"
$ps.Range("D9").Value = 400

$ps.Range("A10").Value = 3
$ps.Range("B10").Value = "Bhira -Mah."
$ps.Range("C10").Value = "1x150"
$ps.Range("D10").Value = 150

$ps.Range("A11").Value = 4
$ps.Range("B11").Value = "Srisailam LBPH -Telangana"
$ps.Range("C11").Value = "6x150"
$ps.Range("D11").Value = 900

$ps.Range("A12").Value = 5
$ps.Range("B12").Value = "Purlia PSS -W.B."
$ps.Range("C12").Value = "4x225"
$ps.Range("D12").Value = 900

$ps.Range("A13").Value = 6
$ps.Range("B13").Value = "Ghatgar -Mah."
$ps.Range("C13").Value = "2x125"
$ps.Range("D13").Value = 250

$ps.Range("B9:C13").Font.Name = "Arial"
$ps.Range("B9:C13").Font.Size = 10

# ---- Working total row 14 --------------------------------------------------
$ps.Range("C14").Value = "Working tot."
$ps.Range("D14").Formula = "=SUM(D8:D13)"

# ---- Section header row 16 -------------------------------------------------
$ps.Range("A16").Value = "A. Schemes under construction"
$ps.Range("A16").Font.Name = "Arial"
$ps.Range("A16").Font.Size = 12
$ps.Range("A16").Font.Bold = $true

# ---- Data rows 17-19 (under construction) ----------------------------------
$ps.Range("A17").Value = 1
$ps.Range("B17").Value = "Tehri St.-II -Uttarakhand "
$ps.Range("C17").Value = "4x250"
$ps.Range("D17").Value = 1000
$ps.Range("E17").Value = "Likely commissioning by 2021-23 (June’22)"
$ps.Range("E17").Font.Name = "Arial"
$ps.Range("E17").Font.Size = 13

$ps.Range("A18").Value = 2
$ps.Range("B18").Value = "Koyna Left Bank -Mah."
$ps.Range("C18").Value = "2x40"
$ps.Range("D18").Value = 80
$ps.Range("E18").Value = "Likely commissioning by 2022-23"

$ps.Range("A19").Value = 3
$ps.Range("B19").Value = "Kundah Pump Storage (Stage I,II,II&IV)-T. N"
$ps.Range("C19").Value = "4x 125"
$ps.Range("D19").Value = 500
$ps.Range("E19").Value = "Likely commissioning by 2022-23"

$ps.Range("B17:C19").Font.Name = "Arial"
$ps.Range("B17:C19").Font.Size = 10

$ps.Range("A4").Select()

# ===========================================================================
# 2. Update the "About" sheet: notes + sources for pumped hydro / geothermal
# ===========================================================================
$about = $wb.Worksheets.Item("About")

# ---- New source column D rows 10-14 (Pumped hydro source) ------------------
$about.Range("D10").Value = "Pumped hydro"
$about.Range("D11").Value = "Central Electricity Authority"
$about.Range("D12").Value = 43709
$about.Range("D12").NumberFormat = "mmm-yy"
$about.Range("D13").Value = "Pumped Storage Development in India"
$about.Range("D14").Value = "http://www.cea.nic.in/reports/monthly/hydro/2019/pump_storage-09.pdf"
$about.Hyperlinks.Add($about.Range("D14"), "http://www.cea.nic.in/reports/monthly/hydro/2019/pump_storage-09.pdf") | Out-Null

# ---- Notes section rewrite (rows 18-26, 28-29) -----------------------------
$about.Range("A18").Value = "In the India EPS, the geothermal plant type is repurposed as pumped hydro capacity."
$about.Range("A19").Value = "Natural gas open cycle plants and pumped hydro are included as peaking plants ."
$about.Range("A20").Value = "No solar thermal or offshore wind capacity existed in 2017."
$about.Range("A21").Value = "No lignite plants were operational in 2017."
$about.Range("A22").Value = "It is assumed that all peaking plants provide flexibility points."
$about.Range("A23").Value = "Diesel is accounted in petroleum subscript of EPS."
$about.Range("A24").Value = "Waste to energy is accounted under Municipal Solid Waste."
$about.Range("A25").Value = "Hydro includes large, small, and pumped storage."
$about.Range("A26").Value = "As per CEA's NEP, of the monitored natural gas based capacity in March 2017, 350 MW is open cycle `nwhich is suited for peaking. "
$about.Rows("26").RowHeight = 24

$about.Range("A27").ClearContents()
$about.Range("A28").Value = "As official data has no detailed breakdown of liquid fuel plants, "
$about.Range("A29").Value = "Crude Oil and Residual Fuel Oil types are assumed to be accounted for in Diesel subscript."
$about.Range("A30").ClearContents()
$about.Range("A31").ClearContents()

# ===========================================================================
# 3. SYC-SYEGC: geothermal capacity now pulls from the Pumped storage sheet
# ===========================================================================
$syegc = $wb.Worksheets.Item("SYC-SYEGC")
$syegc.Range("B10").Formula = "='Pumped storage'!D14"

Write-Host "edit complete"
